# docx_formatter: split each "Qn.  <rest>" run into three separate runs
#   "Qn" | "." | " <rest>"
# and (for Q3) drop the proofErr-wrapped spell-check markup around
# "PortFast", folding it back into plain run text.

function Split-QuestionRun {
    param(
        $Paragraph,
        [string]$Number,
        [string]$Rest
    )

    $range = $Paragraph.Range

    # Recover the paragraph's own opening tag (with its w14:paraId / rsid
    # attributes, etc.) so we don't clobber paragraph identity when we
    # rewrite its contents.
    $openXml = $range.WordOpenXML
    if ($openXml -match '(<w:p\b[^>]*>)') {
        $pOpenTag = $matches[1]
    } else {
        $pOpenTag = '<w:p>'
    }

    $escRest = $Rest -replace '&','&amp;' -replace '<','&lt;' -replace '>','&gt;'

    $xml = $pOpenTag `
        + '<w:r><w:t>Q' + $Number + '</w:t></w:r>' `
        + '<w:r><w:t>.</w:t></w:r>' `
        + '<w:r><w:t xml:space="preserve">' + $escRest + '</w:t></w:r>' `
        + '</w:p>'

    $null = $range.InsertXML($xml)
}

$d = $word.ActiveDocument

Split-QuestionRun $d.Paragraphs(1) "1" " Why can a frame loop endlessly in a switched environment with more than one way to a destination? (if spanning tree or a similar protocol is not enabled)"
Split-QuestionRun $d.Paragraphs(2) "2" " Are there any concerns with allowing auto-selection of the root bridge based on MAC addresses?"
Split-QuestionRun $d.Paragraphs(3) "3" " Please provide the shortest summary of PortFast you can."
Split-QuestionRun $d.Paragraphs(4) "4" " Please provide the shortest summary of BPDU Guard you can."
